$d = $word.ActiveDocument

# Locate the paragraph that ends the "KEY ACHIEVEMENTS AND IMPACT" bullet
# list ("Expert methodology validated at highest judicial level") so the
# two new bullet paragraphs can be inserted immediately after it and
# before the following "TECHNICAL SKILLS" heading.
$anchorIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "*Expert methodology validated at highest judicial level*") {
        $anchorIdx = $i
        break
    }
}

# Insert a new empty paragraph right after the anchor paragraph and fill
# it with the first new bullet (plain text, no special formatting).
$anchorRange = $d.Paragraphs.Item($anchorIdx).Range
$anchorRange.Collapse(0)
$anchorRange.InsertParagraphAfter()
$bullet1Idx = $anchorIdx + 1
$d.Paragraphs.Item($bullet1Idx).Range.Text = "• Breakthrough demographic discovery: Uncovered systematic voter miscoding affecting millions"

# Insert a second new empty paragraph after the first new bullet and fill
# it with the second bullet's full text; "178%" will be bolded/colored
# afterwards.
$bullet1Range = $d.Paragraphs.Item($bullet1Idx).Range
$bullet1Range.Collapse(0)
$bullet1Range.InsertParagraphAfter()
$bullet2Idx = $bullet1Idx + 1
$d.Paragraphs.Item($bullet2Idx).Range.Text = "• 178% accuracy improvement in racial classification algorithms"

# Bold + color just the "178%" run, matching the rest of the document's
# emphasized-statistic style (bold, RGB 2C3E50).
$findRange = $d.Paragraphs.Item($bullet2Idx).Range.Duplicate
$found = $findRange.Find.Execute("178%", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $findRange.Font.Bold = 1
    $findRange.Font.Color = 5258796
}
